$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheets
$ws1.Name = "Data"
$ws2.Name = "Metadata"

# --- Update "Data" sheet rows 2-19 (years 2024 down to 2006) ---
$ws1.Cells.Item(2,1).Value = '2024'
$ws1.Cells.Item(2,2).Value = 32
$ws1.Cells.Item(2,3).Value = 59.4
$ws1.Cells.Item(2,4).Value = 3.1
$ws1.Cells.Item(2,5).Value = 4.5
$ws1.Cells.Item(2,6).Value = 0.053
$ws1.Cells.Item(2,7).Value = 3.6
$ws1.Cells.Item(2,8).Value = 0.99

$ws1.Cells.Item(3,1).Value = '2023'
$ws1.Cells.Item(3,2).Value = 32.5
$ws1.Cells.Item(3,3).Value = 59.7
$ws1.Cells.Item(3,4).Value = 2.8
$ws1.Cells.Item(3,5).Value = 4
$ws1.Cells.Item(3,6).Value = 0.022
$ws1.Cells.Item(3,7).Value = 3.2
$ws1.Cells.Item(3,8).Value = 0.98

$ws1.Cells.Item(4,1).Value = '2022'
$ws1.Cells.Item(4,2).Value = 31.1
$ws1.Cells.Item(4,3).Value = 60.1
$ws1.Cells.Item(4,4).Value = 3.2
$ws1.Cells.Item(4,5).Value = 4.3
$ws1.Cells.Item(4,6).Value = 0.028
$ws1.Cells.Item(4,7).Value = 3.2
$ws1.Cells.Item(4,8).Value = 1.2

$ws1.Cells.Item(5,1).Value = '2021'
$ws1.Cells.Item(5,2).Value = 31.4
$ws1.Cells.Item(5,3).Value = 57.6
$ws1.Cells.Item(5,4).Value = 2.4
$ws1.Cells.Item(5,5).Value = 3.8
$ws1.Cells.Item(5,6).Value = 0.016
$ws1.Cells.Item(5,7).Value = 3.3
$ws1.Cells.Item(5,8).Value = 1.5

$ws1.Cells.Item(6,1).Value = '2019'
$ws1.Cells.Item(6,2).Value = 30.1
$ws1.Cells.Item(6,3).Value = 57.1
$ws1.Cells.Item(6,4).Value = 2.3
$ws1.Cells.Item(6,5).Value = 3.9
$ws1.Cells.Item(6,6).Value = 0.2
$ws1.Cells.Item(6,7).Value = 5.2
$ws1.Cells.Item(6,8).Value = 1.2

$ws1.Cells.Item(7,1).Value = '2018'
$ws1.Cells.Item(7,2).Value = 29.2
$ws1.Cells.Item(7,3).Value = 57.9
$ws1.Cells.Item(7,4).Value = 2.4
$ws1.Cells.Item(7,5).Value = 4
$ws1.Cells.Item(7,6).Value = 0.18
$ws1.Cells.Item(7,7).Value = 4.9
$ws1.Cells.Item(7,8).Value = 1.3

$ws1.Cells.Item(8,1).Value = '2017'
$ws1.Cells.Item(8,2).Value = 28.4
$ws1.Cells.Item(8,3).Value = 58.4
$ws1.Cells.Item(8,4).Value = 2.3
$ws1.Cells.Item(8,5).Value = 4.1
$ws1.Cells.Item(8,6).Value = 0.15
$ws1.Cells.Item(8,7).Value = 5.1
$ws1.Cells.Item(8,8).Value = 1.4

$ws1.Cells.Item(9,1).Value = '2016'
$ws1.Cells.Item(9,2).Value = 28.3
$ws1.Cells.Item(9,3).Value = 58
$ws1.Cells.Item(9,4).Value = 2.1
$ws1.Cells.Item(9,5).Value = 4.5
$ws1.Cells.Item(9,6).Value = 0.22
$ws1.Cells.Item(9,7).Value = 5.4
$ws1.Cells.Item(9,8).Value = 1.6

$ws1.Cells.Item(10,1).Value = '2015'
$ws1.Cells.Item(10,2).Value = 27.1
$ws1.Cells.Item(10,3).Value = 58.4
$ws1.Cells.Item(10,4).Value = 2.1
$ws1.Cells.Item(10,5).Value = 4.4
$ws1.Cells.Item(10,6).Value = 0.32
$ws1.Cells.Item(10,7).Value = 6.1
$ws1.Cells.Item(10,8).Value = 1.6

$ws1.Cells.Item(11,1).Value = '2014'
$ws1.Cells.Item(11,2).Value = 27.6
$ws1.Cells.Item(11,3).Value = 57.7
$ws1.Cells.Item(11,4).Value = 2.1
$ws1.Cells.Item(11,5).Value = 4.3
$ws1.Cells.Item(11,6).Value = 0.27
$ws1.Cells.Item(11,7).Value = 6.4
$ws1.Cells.Item(11,8).Value = 1.7

$ws1.Cells.Item(12,1).Value = '2013'
$ws1.Cells.Item(12,2).Value = 27.8
$ws1.Cells.Item(12,3).Value = 57
$ws1.Cells.Item(12,4).Value = 1.9
$ws1.Cells.Item(12,5).Value = 4.6
$ws1.Cells.Item(12,6).Value = 0.29
$ws1.Cells.Item(12,7).Value = 6.4
$ws1.Cells.Item(12,8).Value = 2

$ws1.Cells.Item(13,1).Value = '2012'
$ws1.Cells.Item(13,2).Value = 27.1
$ws1.Cells.Item(13,3).Value = 56.1
$ws1.Cells.Item(13,4).Value = 1.7
$ws1.Cells.Item(13,5).Value = 4.7
$ws1.Cells.Item(13,6).Value = 0.43
$ws1.Cells.Item(13,7).Value = 7.3
$ws1.Cells.Item(13,8).Value = 2.5

$ws1.Cells.Item(14,1).Value = '2011'
$ws1.Cells.Item(14,2).Value = 28.3
$ws1.Cells.Item(14,3).Value = 54.8
$ws1.Cells.Item(14,4).Value = 2
$ws1.Cells.Item(14,5).Value = 4.7
$ws1.Cells.Item(14,6).Value = 0.47
$ws1.Cells.Item(14,7).Value = 7.1
$ws1.Cells.Item(14,8).Value = 2.6

$ws1.Cells.Item(15,1).Value = '2010'
$ws1.Cells.Item(15,2).Value = 31.9
$ws1.Cells.Item(15,3).Value = 52.1
$ws1.Cells.Item(15,4).Value = 1.7
$ws1.Cells.Item(15,5).Value = 5.1
$ws1.Cells.Item(15,6).Value = 0.61
$ws1.Cells.Item(15,7).Value = 5.8
$ws1.Cells.Item(15,8).Value = 2.8

$ws1.Cells.Item(16,1).Value = '2009'
$ws1.Cells.Item(16,2).Value = 30.2
$ws1.Cells.Item(16,3).Value = 52.5
$ws1.Cells.Item(16,4).Value = 1.8
$ws1.Cells.Item(16,5).Value = 5.4
$ws1.Cells.Item(16,6).Value = 0.8
$ws1.Cells.Item(16,7).Value = 6
$ws1.Cells.Item(16,8).Value = 3.3

$ws1.Cells.Item(17,1).Value = '2008'
$ws1.Cells.Item(17,2).Value = 29.2
$ws1.Cells.Item(17,3).Value = 49.6
$ws1.Cells.Item(17,4).Value = 2.3
$ws1.Cells.Item(17,5).Value = 5.9
$ws1.Cells.Item(17,6).Value = 0.63
$ws1.Cells.Item(17,7).Value = 9
$ws1.Cells.Item(17,8).Value = 3.5

$ws1.Cells.Item(18,1).Value = '2007'
$ws1.Cells.Item(18,2).Value = 35.6
$ws1.Cells.Item(18,3).Value = 42.3
$ws1.Cells.Item(18,4).Value = 2.1
$ws1.Cells.Item(18,5).Value = 6
$ws1.Cells.Item(18,6).Value = 1.1
$ws1.Cells.Item(18,7).Value = 9
$ws1.Cells.Item(18,8).Value = 4

$ws1.Cells.Item(19,1).Value = '2006'
$ws1.Cells.Item(19,2).Value = 36.7
$ws1.Cells.Item(19,3).Value = 42.1
$ws1.Cells.Item(19,4).Value = 1.6
$ws1.Cells.Item(19,5).Value = 6.1
$ws1.Cells.Item(19,6).Value = 1.2
$ws1.Cells.Item(19,7).Value = 8.5
$ws1.Cells.Item(19,8).Value = 3.8

# --- Update "Metadata" sheet rows 1-11 ---
$ws2.Cells.Item(1,1).Value = ' '
$ws2.Cells.Item(1,2).Value = ' '

$ws2.Cells.Item(2,1).Value = 'nomindicador'
$ws2.Cells.Item(2,2).Value = 'Distribución porcentual de personas según institución prestadora en la cual declaran tener cobertura vigente'

$ws2.Cells.Item(3,1).Value = 'derecho'
$ws2.Cells.Item(3,2).Value = 'Salud'

$ws2.Cells.Item(4,1).Value = 'conindicador'
$ws2.Cells.Item(4,2).Value = 'Cobertura integral de salud'

$ws2.Cells.Item(5,1).Value = 'tipoind'
$ws2.Cells.Item(5,2).Value = 'Resultados'

$ws2.Cells.Item(6,1).Value = 'definicion'
$ws2.Cells.Item(6,2).Value = 'El indicador refleja la distribución porcentual de personas según institución prestadora integral de salud en la cual declaran tener cobertura (derecho vigente según la pregunta específica de la ECH).'

$ws2.Cells.Item(7,1).Value = 'calculo'
$ws2.Cells.Item(7,2).Value = 'Para cada año calcular: (Cantidad de personas según tipo de institución prestadora de salud en la que tienen derechos vigentes / Cantidad de habitantes)*100'

$ws2.Cells.Item(8,1).Value = 'observaciones'
$ws2.Cells.Item(8,2).Value = 'Desde marzo de 2020 hasta junio de 2021 se interrumpió el relevamiento presencial y se aplicó de manera telefónica un cuestionario restringido con el objetivo de continuar publicando los indicadores de ingresos y mercado de trabajo. En ese período la encuesta pasó a ser de paneles rotativos elegidos al azar a partir de los casos respondentes del año anterior. 
En julio de 2021 el INE retomó la realización de encuestas presenciales, pero introdujo un cambio metodológico, ya que la ECH pasa a ser una encuesta de panel rotativo con periodicidad mensual compuesta por seis paneles o grupos de rotación, cada uno de los cuales es una muestra representativa de la población. Con esta nueva metodología, cada hogar seleccionado participa durante seis meses de la ECH.  
A partir de 2020 cambia el modo de relevar cobertura de salud. Antes de esta fecha se les consultaba a los/as encuestados por cobertura en cada uno de los prestadores posibles. Durante el 2020 y el primer semestre de 2021, se relevó únicamente el principal prestador de salud. En el segundo semestre de 2021 se relevó el prestador principal y secundario, hecho que habilita reconstruir un indicador más próximo al calculado antes de 2019. 
Para 2021, este indicador se calcula únicamente a partir de la implantación de modalidad panel del segundo semestre de 2021. Dados los cambios metodológicos en la formulación de las preguntas, no se incorpora a la serie el año 2020 y tampoco se considera la información del primer semestre de 2021.'

$ws2.Cells.Item(9,1).Value = 'actualizacion'
$ws2.Cells.Item(9,2).Value = 'Julio 2025'

$ws2.Cells.Item(10,1).Value = 'cita'
$ws2.Cells.Item(10,2).Value = 'UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE (Hasta 2019) / A partir de 2020 con base en ECH - INE
'

$ws2.Cells.Item(11,1).Value = 'Mirador DESCA - UMAD/FCS – INDDHH'
$ws2.Cells.Item(11,2).Value = ' '

Write-Output "done"